$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - column F holds "想去人数" (want-to-go count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 625
$ws1.Range("F4").Value = 623
$ws1.Range("F7").Value = 2748
$ws1.Range("F8").Value = 465
$ws1.Range("F9").Value = 7678
$ws1.Range("F10").Value = 199
$ws1.Range("F12").Value = 36
$ws1.Range("F13").Value = 308
$ws1.Range("F14").Value = 44

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 14
$ws2.Range("F4").Value = 3

# Sheet "全部类型" (all types, combined view)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 625
$ws4.Range("F4").Value = 623
$ws4.Range("F7").Value = 14
$ws4.Range("F9").Value = 2748
$ws4.Range("F10").Value = 465
$ws4.Range("F11").Value = 7678
$ws4.Range("F12").Value = 199
$ws4.Range("F14").Value = 36
$ws4.Range("F15").Value = 3
$ws4.Range("F17").Value = 308
$ws4.Range("F18").Value = 44
